$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("A1").Value = "Block"
$ws.Range("B1").Value = "Meso_OS_BAsqft"
$ws.Range("C1").Value = "Category"

# --- Refresh the sort-state bookkeeping: the data is now (re)sorted
#     ascending by the Block column instead of descending by Meso_OS_BAsqft. ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A13")) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:C13"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Rewrite data rows in the final, ascending-by-Block order ---
$data = @(
    @(1,  87,                  "high"),
    @(2,  66.400000000000006,  "med"),
    @(3,  81.900000000000006,  "high"),
    @(4,  52.2,                "low"),
    @(5,  47.9,                "low"),
    @(6,  41.1,                "low"),
    @(7,  61,                  "med"),
    @(8,  64.400000000000006,  "med"),
    @(9,  83.6,                "high"),
    @(10, 64.900000000000006,  "med"),
    @(11, 47,                  "low"),
    @(12, 84.3,                "high")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Remove the old "RAIN LOGGER" marker column D entirely ---
$ws.Range("D1:D13").ClearContents()

# --- Column B width (bestfit-style custom width) ---
$ws.Columns.Item(2).ColumnWidth = 15.14

# --- Update selection to match the new active cell/sqref ---
$ws.Range("A2:C13").Select()
